$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'47.159.35"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "'2.489.72"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'322.01"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'108.54"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "'38.85"
$ws.Range("E10").Value = "  +7.48%  "
$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "'18.35"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "'2.877.22"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "'2.491.87"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "'0.850"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "'47.060.67"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "'12.72"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  +14.97%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'246.27"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D27").Value = "'25.87"
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").Value = "'10.03"
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("D30").Value = "'0.140"
$ws.Range("E30").Value = "  +8.41%  "
$ws.Range("D31").Value = "'35.20"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("D32").Value = "'49.92"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").Value = "'20.05"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("E34").Value = "  +1.67%  "
$ws.Range("E35").Value = "  +2.32%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'4.69"
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.96"
$ws.Range("E38").Value = "  +2.86%  "
$ws.Range("D39").Value = "'2.97"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").Value = "'120.80"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("D45").Value = "'1.993.75"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  +2.26%  "
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("D50").Value = "'5.17"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "'56.54"
$ws.Range("E51").Value = "  +3.34%  "